$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.724.93'
$ws.Range('E2').Value = '  -3.83%  '
$ws.Range('D3').Value = '3.315.64'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '573.22'
$ws.Range('E5').Value = '  -3.05%  '
$ws.Range('D6').Value = '182.74'
$ws.Range('E6').Value = '  -5.40%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('D9').Value = '0.129'
$ws.Range('E9').Value = '  -3.14%  '
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('E11').Value = '  -4.52%  '
$ws.Range('D12').Value = '3.894.58'
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('D13').Value = '0.137'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('D14').Value = '27.21'
$ws.Range('E14').Value = '  -3.42%  '
$ws.Range('D15').Value = '66.758.96'
$ws.Range('E15').Value = '  -3.74%  '
$ws.Range('E16').Value = '  -2.41%  '
$ws.Range('D17').Value = '3.302.05'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').Value = '437.65'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '13.76'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = '5.68'
$ws.Range('E20').Value = '  -2.41%  '
$ws.Range('D21').Value = '7.64'
$ws.Range('E21').Value = '  -1.30%  '
$ws.Range('D22').Value = '73.82'
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('E24').Value = '  -0.28%  '
$ws.Range('D25').Value = '0.0000118'
$ws.Range('E25').Value = '  -2.23%  '
$ws.Range('D26').Value = '0.194'
$ws.Range('E26').Value = '  +1.18%  '
$ws.Range('D27').Value = '9.10'
$ws.Range('E27').Value = '  -5.14%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('D30').Value = '22.86'
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('D31').Value = '5.33'
$ws.Range('E31').Value = '  -4.33%  '
$ws.Range('D33').Value = '6.79'
$ws.Range('E33').Value = '  -2.89%  '
$ws.Range('E34').Value = '  -3.72%  '
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('D36').Value = '160.37'
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').Value = '1.86'
$ws.Range('E37').Value = '  -3.14%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '27.39'
$ws.Range('E38').Value = '  +1.37%  '
$ws.Range('D39').Value = '2.812.41'
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('D40').Value = '0.792'
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('D41').Value = '4.46'
$ws.Range('E41').Value = '  -2.51%  '
$ws.Range('E42').Value = '  -3.67%  '
$ws.Range('D43').Value = '0.0676'
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('D44').Value = '40.17'
$ws.Range('E44').Value = '  -2.41%  '
$ws.Range('D45').Value = '24.35'
$ws.Range('E45').Value = '  -4.04%  '
$ws.Range('E46').Value = '  -6.43%  '
$ws.Range('D47').Value = '318.89'
$ws.Range('E47').Value = '  -7.64%  '
$ws.Range('D48').Value = '0.0273'
$ws.Range('E48').Value = '  -3.28%  '
$ws.Range('D49').Value = '0.987'
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('D51').Value = '0.0997'
$ws.Range('E51').Value = '  -1.62%  '
